# Update scaling mappings to proper format
#
# The "year" mapping sheet (sheet3 / "year") is extended with four new
# columns (select_scaling_year, start_scaling_year, end_scaling_year,
# Comment) and its example/default row is updated to use real sample
# values instead of placeholder "NA" text. The "year" sheet also becomes
# the active / selected tab of the workbook (it previously was the
# "map" sheet).

$wb = $excel.ActiveWorkbook

$mapSheet  = $wb.Worksheets.Item(1)   # "map"
$yearSheet = $wb.Worksheets.Item(3)   # "year"

# --- Add the new header columns (row 1) -----------------------------------
# Must be written before the row-2 values below so that new shared-string
# entries are appended in header-then-data order.
$yearSheet.Range("E1").Value = "select_scaling_year"
$yearSheet.Range("F1").Value = "start_scaling_year"
$yearSheet.Range("G1").Value = "end_scaling_year"
$yearSheet.Range("H1").Value = "Comment"

# --- Update / add the example row (row 2) ----------------------------------
$yearSheet.Range("A2").Value = "mkd"
$yearSheet.Range("B2").Value = "all"
$yearSheet.Range("E2").Value = "NA"
$yearSheet.Range("F2").Value = 1990
$yearSheet.Range("G2").Value = 2010
$yearSheet.Range("H2").Value = "Don't scale to 1990-1991 drop so as to be closer to EMEP trend"

# --- Make the "year" sheet the active / selected tab ------------------------
$yearSheet.Activate()
$yearSheet.Rows("1:2").EntireRow.Select()
